$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 27.02.2022 06:45"

# Update row 8 (Benzina Albert Modrice) values
$ws.Range("B8").Value = 39.5
$ws.Range("C8").Value = 37.9

# D8/E8 become text in the edited file (rather than numbers), so force
# text storage (otherwise "+1.6" / the datetime string would silently be
# re-interpreted as a number), then drop the temporary Text format again
# so the cells end up unstyled, matching the target.
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "+1.6"
$ws.Range("E8").Value = "2022-02-27 06:45:16"
$ws.Range("D8:E8").ClearFormats()
